# Update cryptos list cell values per the scraped diff (Sat Nov 4 13:58:44 UTC 2023 GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.028.74'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '1.851.10'
$ws.Range("E3").Value = '  +2.00%  '
$ws.Range("E4").Value = '  +0.28%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '237.84'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +3.20%  '
$ws.Range("E6").Value = '  +0.78%  '
$ws.Range("E7").Value = '  +0.19%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '42.18'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +6.04%  '
$ws.Range("E9").Value = '  +1.85%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0692'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.66%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0991'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '2.118.60'
$ws.Range("E12").Value = '  +2.05%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '11.39'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '1.843.47'
$ws.Range("E14").Value = '  +1.80%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.673'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.07%  '
$ws.Range("E16").Value = '  +3.80%  '
$ws.Range("D17").Value = '35.042.12'
$ws.Range("E17").Value = '  +0.80%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '70.11'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("E19").Value = '  +1.65%  '
$ws.Range("E20").Value = '  +0.08%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '12.14'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  +2.01%  '
$ws.Range("E22").Value = '  +2.50%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("E24").Value = '  +1.38%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '169.76'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -2.20%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '7.98'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +3.02%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '1.80'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +18.93%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '17.59'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.68%  '
$ws.Range("E29").Value = '  +0.25%  '
$ws.Range("E30").Value = '  +0.25%  '
$ws.Range("E31").Value = '  +1.95%  '
$ws.Range("E32").Value = '  +0.30%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '4.01'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +2.33%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.67'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +23.46%  '
$ws.Range("E35").Value = '  +10.51%  '
$ws.Range("E36").Value = '  +5.38%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.776'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +12.88%  '
$ws.Range("E38").Value = '  +10.09%  '
$ws.Range("E39").Value = '  +5.33%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '90.10'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("D41").Value = '1.347.06'
$ws.Range("E41").Value = '  +1.28%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '14.70'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +3.48%  '
$ws.Range("E43").Value = '  +2.69%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("B45").Value = 'Kaspa'
$ws.Range("C45").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.0556'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +6.97%  '
$ws.Range("B46").Value = 'MXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.74'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.30%  '
$ws.Range("B47").Value = 'Gas'
$ws.Range("C47").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '12.18'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +43.39%  '
$ws.Range("E48").Value = '  +6.72%  '
$ws.Range("D49").Value = '2.038.52'
$ws.Range("E49").Value = '  +2.23%  '
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0689'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +4.20%  '
$ws.Range("E51").Value = '  +0.24%  '
